$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Main St Radford Virginia")
$ws.Range("E2").Value = 0.081
$ws.Range("E3").Value = 0.081
$ws.Range("E4").Value = 0.081
$ws.Range("G4").Value = 0.0267
$ws.Range("J4").Value = 0.0267
$ws.Range("K4").Value = 0.0135
$ws.Range("L4").Value = 0.0137
$ws.Range("N4").Value = 0.0271
$ws.Range("O4").Value = 0.0274
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("O7").ClearContents()
$ws.Range("E8").Value = 0.0733
$ws.Range("E9").Value = 0.0733
$ws.Range("E10").Value = 0.0733
$ws.Range("G10").Value = 0.0196
$ws.Range("H10").Value = 0.0049
$ws.Range("J10").Value = 0.0293
$ws.Range("K10").Value = 0.0049
$ws.Range("L10").Value = 0.0246
$ws.Range("M10").Value = 0.0148
$ws.Range("N10").Value = 0.0441
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0.010475
$ws.Range("Q10").Value = 0.010475
$ws.Range("R10").Value = 0.031425
$ws.Range("S10").Value = 0.010475
$ws.Range("T10").Value = 0.010475
$ws.Range("U10").Value = 0.010475
$ws.Range("V10").Value = 0.031425
$ws.Range("W10").Value = 0.1257

$ws = $wb.Worksheets.Item("Marengo Illinois")
$ws.Range("O7").ClearContents()

$ws = $wb.Worksheets.Item("Modrice Czech Republic")
$ws.Range("E2").Value = 0.0446
$ws.Range("E3").Value = 0.0446
$ws.Range("E4").Value = 0.0446
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("O7").ClearContents()
$ws.Range("E8").Value = 0.0114
$ws.Range("E9").Value = 0.0114
$ws.Range("E10").Value = 0.0114
$ws.Range("O10").Value = 0.0057
$ws.Range("P10").Value = 0.001625
$ws.Range("Q10").Value = 0.001625
$ws.Range("R10").Value = 0.004875
$ws.Range("S10").Value = 0.001625
$ws.Range("T10").Value = 0.001625
$ws.Range("U10").Value = 0.001625
$ws.Range("V10").Value = 0.004875
$ws.Range("W10").Value = 0.0195

$ws = $wb.Worksheets.Item("Nova Milanese, Italy")
$ws.Range("E2").Value = 0.1
$ws.Range("E3").Value = 0.1
$ws.Range("E4").Value = 0.1
$ws.Range("O4").Value = 0.05
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0

$ws = $wb.Worksheets.Item("Ratingen Germany")
$ws.Range("E2").Value = 0.0338
$ws.Range("E3").Value = 0.0338
$ws.Range("E4").Value = 0.0338
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0

$ws = $wb.Worksheets.Item("Rock Road (KCS) Radford Virgin")
$ws.Range("E2").Value = 0.1754
$ws.Range("E3").Value = 0.1754
$ws.Range("E4").Value = 0.1754
$ws.Range("K4").Value = 0.029
$ws.Range("L4").Value = 0.0455
$ws.Range("M4").Value = 0.0152
$ws.Range("N4").Value = 0.0896
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("O7").ClearContents()

$ws = $wb.Worksheets.Item("Rock Road Radford Virginia")
$ws.Range("E2").Value = 0.0612
$ws.Range("E3").Value = 0.0612
$ws.Range("E4").Value = 0.0612
$ws.Range("I4").Value = 0.0087
$ws.Range("J4").Value = 0.0087
$ws.Range("L4").Value = 0.0268
$ws.Range("M4").Value = 0.0268
$ws.Range("N4").Value = 0.0531
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("E5").Value = 0.818181818181818
$ws.Range("E6").Value = 0.818181818181818
$ws.Range("E7").Value = 0.818181818181818
$ws.Range("O7").Value = 1
$ws.Range("P7").Value = 0.818181818181818
$ws.Range("Q7").Value = 0.818181818181818
$ws.Range("R7").Value = 0.818181818181818
$ws.Range("S7").Value = 0.818181818181818
$ws.Range("T7").Value = 0.818181818181818
$ws.Range("U7").Value = 0.818181818181818
$ws.Range("V7").Value = 0.818181818181818
$ws.Range("W7").Value = 0.818181818181818
$ws.Range("E8").Value = 0.0861
$ws.Range("E9").Value = 0.0861
$ws.Range("E10").Value = 0.0861
$ws.Range("G10").Value = 0.0301
$ws.Range("I10").Value = 0.0074
$ws.Range("J10").Value = 0.0373
$ws.Range("L10").Value = 0.021
$ws.Range("M10").Value = 0.014
$ws.Range("N10").Value = 0.0352
$ws.Range("O10").Value = 0.0135
$ws.Range("P10").Value = 0.0123
$ws.Range("Q10").Value = 0.0123
$ws.Range("R10").Value = 0.0369
$ws.Range("S10").Value = 0.0123
$ws.Range("T10").Value = 0.0123
$ws.Range("U10").Value = 0.0123
$ws.Range("V10").Value = 0.0369
$ws.Range("W10").Value = 0.1476

$ws = $wb.Worksheets.Item("SEEPZ-SEZ Mumbai India")
$ws.Range("E2").Value = 0.7407
$ws.Range("E3").Value = 0.7407
$ws.Range("E4").Value = 0.7407
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.0833333333333333
$ws.Range("Q4").Value = 0.0833333333333333
$ws.Range("R4").Value = 0.25
$ws.Range("S4").Value = 0.0833333333333333
$ws.Range("T4").Value = 0.0833333333333333
$ws.Range("U4").Value = 0.0833333333333333
$ws.Range("V4").Value = 0.25
$ws.Range("W4").Value = 1

$ws = $wb.Worksheets.Item("Santa Barbara California")
$ws.Range("E2").Value = 0.1258
$ws.Range("E3").Value = 0.1258
$ws.Range("E4").Value = 0.1258
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0

$ws = $wb.Worksheets.Item("Shanghai Minhang District Chin")
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0

$ws = $wb.Worksheets.Item("Tianjin China")
$ws.Range("E2").Value = 0.0113
$ws.Range("E3").Value = 0.0113
$ws.Range("E4").Value = 0.0113
$ws.Range("O4").Value = 0.0114
$ws.Range("O7").ClearContents()
$ws.Range("E8").Value = 0.0249
$ws.Range("E9").Value = 0.0249
$ws.Range("E10").Value = 0.0249
$ws.Range("H10").Value = 0.025
$ws.Range("J10").Value = 0.0248
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0.00355833333333333
$ws.Range("Q10").Value = 0.00355833333333333
$ws.Range("R10").Value = 0.010675
$ws.Range("S10").Value = 0.00355833333333333
$ws.Range("T10").Value = 0.00355833333333333
$ws.Range("U10").Value = 0.00355833333333333
$ws.Range("V10").Value = 0.010675
$ws.Range("W10").Value = 0.0427

$ws = $wb.Worksheets.Item("Istanbul Turkey")
$ws.Range("E2").Value = 0.4167
$ws.Range("E3").Value = 0.4167
$ws.Range("E4").Value = 0.4167
$ws.Range("J4").Value = 0.1156
$ws.Range("O4").Value = 0.1
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("O5").ClearContents()

$ws = $wb.Worksheets.Item("Khed Taluka India")
$ws.Range("O5").ClearContents()

$ws = $wb.Worksheets.Item("Kongegårdsgatan Molndal Sweden")
$ws.Range("E2").Value = 0.1389
$ws.Range("E3").Value = 0.1389
$ws.Range("E4").Value = 0.1389
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("O7").ClearContents()
